$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.985095739364624
$ws.Range("B1").Value = 2.205040693283081
$ws.Range("C1").Value = 5.043397426605225
$ws.Range("D1").Value = 1.796114206314087
$ws.Range("E1").Value = 1.298387765884399
